$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new column P (16) of data for year 2022, mirroring the existing
# year-over-year columns (D..O) already present in the table.

# Row 3: year header
$ws.Range("P3").Value = 2022

# Row 4: count of reporting insurance companies
$ws.Range("P4").Value = 15

# Row 5: insurance premiums (mln. soms)
$ws.Range("P5").Value = 2130.4

# Copy formatting from the preceding column (O) so the new column matches
# the existing style used throughout the table.
$ws.Range("O3").Copy() | Out-Null
$ws.Range("P3").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("O4").Copy() | Out-Null
$ws.Range("P4").PasteSpecial(-4122) | Out-Null

$ws.Range("O5").Copy() | Out-Null
$ws.Range("P5").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# Update the active selection to match the post-edit state.
$ws.Range("P6").Select() | Out-Null
